$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "M2DocEvaluator.caseConditional(M2DocEvaluator.java:1313)" "M2DocEvaluator.caseConditional(M2DocEvaluator.java:1318)"
Replace-Text "M2DocEvaluator.doSwitch(M2DocEvaluator.java:1084)" "M2DocEvaluator.doSwitch(M2DocEvaluator.java:1096)"
Replace-Text "M2DocEvaluator.caseBlock(M2DocEvaluator.java:1300)" "M2DocEvaluator.caseBlock(M2DocEvaluator.java:1305)"
Replace-Text "M2DocEvaluator.caseConditional(M2DocEvaluator.java:1324)" "M2DocEvaluator.caseConditional(M2DocEvaluator.java:1329)"
Replace-Text "M2DocEvaluator.caseDocumentTemplate(M2DocEvaluator.java:278)" "M2DocEvaluator.caseDocumentTemplate(M2DocEvaluator.java:283)"
Replace-Text "M2DocEvaluator.generate(M2DocEvaluator.java:267)" "M2DocEvaluator.generate(M2DocEvaluator.java:272)"
Replace-Text "AbstractTemplatesTestSuite.prepareoutputAndGenerate(AbstractTemplatesTestSuite.java:475)" "AbstractTemplatesTestSuite.prepareoutputAndGenerate(AbstractTemplatesTestSuite.java:479)"
Replace-Text "AbstractTemplatesTestSuite.generation(AbstractTemplatesTestSuite.java:384)" "AbstractTemplatesTestSuite.generation(AbstractTemplatesTestSuite.java:388)"

Write-Output "Done"
